# Update betting odds values on Sheet1 to reflect the latest FlashScore
# snapshot (rows 2, 7, 8 and 14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Novorizontino - Sport Recife)
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 3.75
$ws.Range("J2").Value = 3.2
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("X2").Value = 9
$ws.Range("AW2").Value = 5.5

# Row 7 (Tepatitlan de Morelos - Tampico Madero)
$ws.Range("H7").Value = 2.82
$ws.Range("I7").Value = 3.3
$ws.Range("L7").Value = 3.75
$ws.Range("N7").Value = 6.65
$ws.Range("P7").Value = 2.52
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.55
$ws.Range("S7").Value = 1.45
$ws.Range("U7").Value = 1.82
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 6.4
$ws.Range("X7").Value = 10.5
$ws.Range("AB7").Value = 35
$ws.Range("AF7").Value = 75
$ws.Range("AG7").Value = 8.75
$ws.Range("AH7").Value = 17.5
$ws.Range("AJ7").Value = 50
$ws.Range("AK7").Value = 32
$ws.Range("AR7").Value = 90
$ws.Range("AT7").Value = 2.35
$ws.Range("AU7").Value = 6.7
$ws.Range("AW7").Value = 5.1
$ws.Range("AX7").Value = 18

# Row 8 (Tacuary - Sp. Luqueno)
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.33

# Row 14 (Zamora - Monagas)
$ws.Range("G14").Value = 2.95
$ws.Range("H14").Value = 3
$ws.Range("J14").Value = 3.45
$ws.Range("K14").Value = 2.02
$ws.Range("P14").Value = 2.65
$ws.Range("S14").Value = 1.4
$ws.Range("T14").Value = 2.5
$ws.Range("U14").Value = 1.78
$ws.Range("V14").Value = 1.83
$ws.Range("W14").Value = 8.25
$ws.Range("X14").Value = 15
$ws.Range("Z14").Value = 37
$ws.Range("AA14").Value = 27
$ws.Range("AD14").Value = 5.8
$ws.Range("AG14").Value = 7.1
$ws.Range("AH14").Value = 11.25
$ws.Range("AK14").Value = 22
$ws.Range("AN14").Value = 4.8
$ws.Range("AO14").Value = 16
$ws.Range("AP14").Value = 22
$ws.Range("AQ14").Value = 75
$ws.Range("AU14").Value = 6.7
$ws.Range("AV14").Value = 60
$ws.Range("AX14").Value = 12.5
